$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows after row 100 (this shifts the former rows 101-133 down to 103-135)
$ws.Rows("101:102").Insert()

# New row 101
$ws.Cells.Item(101, 1).Value = 7
$ws.Cells.Item(101, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(101, 3).Value = "Ñuble"
$ws.Cells.Item(101, 4).Value = 44964
$ws.Cells.Item(101, 5).Value = 16
$ws.Cells.Item(101, 6).Value = 100112021
$ws.Cells.Item(101, 7).Value = "Ají"
$ws.Cells.Item(101, 8).Value = "Chilena(o)"
$ws.Cells.Item(101, 9).Value = "Primera"
$ws.Cells.Item(101, 10).Value = 30
$ws.Cells.Item(101, 11).Value = 23000
$ws.Cells.Item(101, 12).Value = 23000
$ws.Cells.Item(101, 13).Value = 23000
$ws.Cells.Item(101, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(101, 15).Value = "Región del Maule"
$ws.Cells.Item(101, 16).Value = 920
$ws.Cells.Item(101, 17).Value = 25
$ws.Cells.Item(101, 18).Value = "Hortaliza"

# New row 102
$ws.Cells.Item(102, 1).Value = 7
$ws.Cells.Item(102, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(102, 3).Value = "Ñuble"
$ws.Cells.Item(102, 4).Value = 44964
$ws.Cells.Item(102, 5).Value = 16
$ws.Cells.Item(102, 6).Value = 100112021
$ws.Cells.Item(102, 7).Value = "Ají"
$ws.Cells.Item(102, 8).Value = "Cristal"
$ws.Cells.Item(102, 9).Value = "Primera"
$ws.Cells.Item(102, 10).Value = 30
$ws.Cells.Item(102, 11).Value = 23000
$ws.Cells.Item(102, 12).Value = 23000
$ws.Cells.Item(102, 13).Value = 23000
$ws.Cells.Item(102, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(102, 15).Value = "Región del Maule"
$ws.Cells.Item(102, 16).Value = 920
$ws.Cells.Item(102, 17).Value = 25
$ws.Cells.Item(102, 18).Value = "Hortaliza"

Write-Output "Rows inserted and populated successfully"
